$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snippets")

# [Word] (Body) Remove duplicate mapping for select
# Row 26 is a duplicate "Body.select" mapping (the real one is row 20).
# Deleting the whole row shifts every row below it up by one, shrinking the
# table/used-range from A1:F203 down to A1:F202.
$ws.Rows.Item(26).Delete()

# Reset the frozen-pane view / active selection back to the top of the data.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A2").Select() | Out-Null
